$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style/formatting of existing
# header cells (e.g. G1) by copying formats over after setting the value.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add values for the new "Save" column in the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
